$d = $word.ActiveDocument

function Replace-Text($old, $new, [bool]$wholeWord = $false) {
    $d.Content.Find.Execute($old, $true, $wholeWord, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---- Title ----
Replace-Text "Unveiling the Secrets of Dark Matter" "Transcending Time: An Exploration of Historical Eras"

# ---- Author ----
Replace-Text "Neil DeGrasse Tyson" "Sophia Hendricks"

# ---- Email address (local@domain.tld split across runs) ----
Replace-Text "neildegrassetyson_astrophysicist@space" "sophia"
Replace-Text "com" "hendricks@school.edu" $true

# ---- Body paragraph 1 (first big paragraph) ----
Replace-Text "Hidden within the vast expanse of the cosmos lies a mysterious entity known as dark matter, an enigmatic substance that has eluded our full understanding" "The tapestry of history is a captivating work of art, intricately woven from the threads of human experience"

Replace-Text " It permeates the universe, exerting a gravitational influence beyond that accounted for by visible matter" " In this realm of exploration, we traverse vast landscapes of time, unveiling the triumphs, tragedies, and pivotal moments that shape our world"

Replace-Text " Despite its invisible nature, dark matter plays a significant role in shaping the structure and evolution of galaxies, clusters, and the cosmos as a whole" " Delving into the annals of bygone eras, we uncover secrets locked away in dusty tomes, whispered tales echoing through the corridors of forgotten civilizations"

Replace-Text " Embarking on a journey to unravel the secrets of dark matter is a quest to comprehend one of the greatest cosmic enigmas of our time" " This journey of discovery ignites our imagination, unveiling the tapestry of stories and events that shape our present. As we immerse ourselves in the lives, struggles, and aspirations of those who came before us, we gain an invaluable perspective on our own existence"

Replace-Text "Through a series of detailed observations and advanced simulations, astronomers and physicists have pieced together clues about the properties and behavior of dark matter" "The chronicle of history showcases the resilience of the human spirit in the face of adversity"

Replace-Text " Studies of galaxy dynamics and gravitational lensing have revealed the presence of a substantial amount of unseen mass, suggesting the existence of dark matter halos surrounding galaxies" " From ancient empires to modern-day struggles, we marvel at the achievements, lament the heartbreaks, and reflect on the choices that have shaped our trajectory"

Replace-Text " Its gravitational effects have been detected in the motion of stars within galaxies, providing evidence for its invisible presence" " By examining the mistakes of the past, we arm ourselves with wisdom, ensuring we do not repeat them"

Replace-Text " Furthermore, observations of the cosmic microwave background radiation, the remnant glow from the early universe, have hinted at the possible nature of dark matter as a cold and collisionless material" " Through the lens of hindsight, we glean valuable lessons, shedding light on the labyrinthine corridors of human ambition, conflict, and resolution. The study of history is not merely an academic pursuit; it is a transformative experience that molds our understanding of the world, inspiring us to make meaningful contributions to society"

# ---- Body paragraph 2 ----
Replace-Text "The nature of dark matter remains a tantalizing mystery, challenging our current understanding of physics" "Our appreciation deepens as we explore the interconnectedness of historical events"

Replace-Text " Some theories propose that it consists of weakly interacting particles, such as axions or sterile neutrinos, while others suggest that it is a manifestation of extra dimensions beyond our observable realm" " The actions, decisions, and circumstances of one era reverberate through time, impacting future generations"

Replace-Text " These exotic candidates for dark matter continue to drive innovation in theoretical physics and experimental techniques, as scientists strive to unravel the elusive secrets of this enigmatic substance" " Like ripples in a pond, historical events create concentric circles of consequences, shaping the cultural, political, and societal landscapes we witness today. Examining these connections provides a profound understanding of the complexities that have led us to our present juncture, enabling us to appreciate the delicate equilibrium that exists between cause and effect. Through this scholarly expedition, we gain a holistic perspective, unraveling the intricate web that weaves the past, present, and future into an indissoluble tapestry"

# ---- Summary body paragraph ----
Replace-Text "The quest to understand dark matter is a captivating journey into the unknown, pushing the boundaries of our knowledge about the universe" "In this voyage through the annals of history, we illuminate the interconnectedness of events, discover the resilience of the human spirit, and gain a profound understanding of the present"

Replace-Text " Through observations, simulations, and theoretical exploration, scientists are inching closer to unraveling the mysteries surrounding this invisible entity" " History provides an invaluable lens through which we explore the triumphs, failures, and aspirations of humanity, gleaning wisdom to ensure we do not repeat the mistakes of the past and igniting the spark of inspiration for a brighter future"

Replace-Text " From galaxy dynamics to the cosmic microwave background radiation, a multitude of clues have emerged, hinting at the existence and properties of dark matter. The hunt for the nature of dark matter continues, promising to reshape our comprehension of the cosmos and deepen our understanding of the fundamental laws governing the universe" " Through this journey of enlightenment, we become custodians of the past, preserving the lessons learned and carrying the torch of knowledge forward into the unfolding chapters of our shared narrative"

# ---- Move the lastRenderedPageBreak marker: it used to sit on the "multitude of clues"
#      run (removed above); the edited document instead carries it on the "Summary" heading run.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $txt = $para.Range.Text
    if ($txt.TrimEnd([char]13) -eq "Summary") {
        $para.Range.Find.Execute("Summary", $true, $false, $false, $false, $false, $true, 1, $false, "Summary", 2) | Out-Null
        break
    }
}

# ---- Append a new empty paragraph at the very end of the body ----
$d.Content.InsertParagraphAfter() | Out-Null
